$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.596.06'
$ws.Range('E2').Value = '  +0.68%  '

$ws.Range('D3').Value = '1.634.90'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.44'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.28%  '

$ws.Range('E6').Value = '  -1.49%  '

$ws.Range('E7').Value = '  +0.15%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.93'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.44%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.257'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.40%  '

$ws.Range('E10').Value = '  -0.06%  '

$ws.Range('E11').Value = '  +0.35%  '

$ws.Range('D12').Value = '1.867.82'
$ws.Range('E12').Value = '  -0.37%  '

$ws.Range('D13').Value = '1.601.67'
$ws.Range('E13').Value = '  -2.35%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.03'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.20%  '

$ws.Range('E15').Value = '  -1.59%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.47'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.13%  '

$ws.Range('D17').Value = '27.605.31'
$ws.Range('E17').Value = '  +0.80%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.95'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.06%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.75'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.33%  '

$ws.Range('E21').Value = '  +0.13%  '

$ws.Range('E22').Value = '  -0.97%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.00'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.51%  '

$ws.Range('E24').Value = '  -2.80%  '

$ws.Range('E25').Value = '  +2.09%  '

$ws.Range('E26').Value = '  -1.06%  '

$ws.Range('E27').Value = '  -1.43%  '

$ws.Range('E28').Value = '  +0.08%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.60'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.14%  '

$ws.Range('E30').Value = '  +0.14%  '

$ws.Range('E32').Value = '  +0.29%  '

$ws.Range('D33').Value = '1.451.64'
$ws.Range('E33').Value = '  +2.58%  '

$ws.Range('E34').Value = '  -1.32%  '

$ws.Range('E35').Value = '  -0.64%  '

$ws.Range('E36').Value = '  +0.02%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.563'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.41%  '

$ws.Range('E38').Value = '  -0.91%  '

$ws.Range('E39').Value = '  +0.44%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.889'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.84%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.82'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.13%  '

$ws.Range('E42').Value = '  +0.11%  '

$ws.Range('E43').Value = '  -0.95%  '

$ws.Range('E44').Value = '  +1.67%  '

$ws.Range('E45').Value = '  +1.10%  '

$ws.Range('E46').Value = '  -0.01%  '

$ws.Range('D47').Value = '1.777.59'
$ws.Range('E47').Value = '  -0.38%  '

$ws.Range('E48').Value = '  +2.13%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '86.21'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.03%  '

$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0984'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.66%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.76'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.46%  '
